$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.203.39"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.479.67"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.482.52"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "4.085.30"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.24%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "67.269.30"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "3.482.32"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.535"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "2.835.56"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0300"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "334.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -2.19%  "
